$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: "Gen" -> "MaxFES"
$ws.Range("A1").Value = "MaxFES"

# Column A values (rows 2-14): change from generation counts to fraction-of-budget values
$colAValues = @(0, 0.001, 0.01, 0.1, 0.2, 0.3, 0.4, 0.5, 0.6, 0.7, 0.8, 0.9, 1)
for ($i = 0; $i -lt $colAValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $colAValues[$i]
}

# Remove the "Run 50" column (originally column AZ) entirely. This shifts the
# old "Mean" column (originally BA) left into AZ, becoming the new last column,
# and drops the now-unused "Run 50" label.
$ws.Range("AZ:AZ").EntireColumn.Delete()

# Recompute the Mean column (now AZ) as the average of the remaining 50 run
# columns (B:AY) for each data row, since Run 50 no longer contributes.
$meanValues = @(48.95120863, 46.71808378, 44.37439163, 41.09028938, 40.46297162, 39.76059126, 39.44682114, 39.32232351, 39.09805572, 38.88135552, 38.80342298, 38.51166529, 38.42848107)
for ($i = 0; $i -lt $meanValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 52).Value = $meanValues[$i]
}
